$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date-like text in column C (and, for symmetry, D) to be
# interpreted as literal text rather than being auto-converted into a
# serial date number when assigned via .Value.
$ws.Range("C3:D3").NumberFormat = "@"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Teixeira"
$ws.Range("C3").Value = "05/06/2024"
$ws.Range("D3").Value = "21/06/2024"
$ws.Range("E3").Value = "nacional"
$ws.Range("F3").Value = "trabalho"
$ws.Range("G3").Value = "f awsfawfawrew"

# Drop the temporary text formatting again so the new row ends up with no
# explicit cell style, matching the source row's plain (unstyled) cells.
$ws.Range("C3:D3").ClearFormats()
